# Add 2022-Q3 data
# -----------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q3" right after "总计" and before
#    "2022-Q2", with the same column layout as the other quarterly sheets.
# 2) Update the "总计" (summary) sheet with a new row for 2022-Q3.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$totalSheetName = $wb.Worksheets.Item(1).Name
$totalSheet = $wb.Worksheets.Item($totalSheetName)

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet, positioned right after "总计"
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# Fetch sheets by NAME (index-based refs shift once a sheet has been
# inserted, so always re-resolve by name after the insert).
$totalSheet = $wb.Worksheets.Item($totalSheetName)
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Copy the column/row formatting of the existing "2022-Q2" sheet (header
# row style + index-column style) onto the new sheet so new cells pick up
# matching styling (border + bold header, etc.).
$q2Sheet.Range("A1:H6").Copy()
$q3Sheet.Range("A1:H6").PasteSpecial(-4122)

# Header row
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Data rows (index column A, then B..H)
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "'118002"
$q3Sheet.Range("C2").Value = "易方达标普全球高端消费品指数增强A（QDII）人民币"
$q3Sheet.Range("D2").Value = "'1.85"
$q3Sheet.Range("E2").Value = "'93.04"
$q3Sheet.Range("F2").Value = "'9.69"
$q3Sheet.Range("G2").Value = "'0.1793"
$q3Sheet.Range("H2").Value = 1

$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "'000593"
$q3Sheet.Range("C3").Value = "易方达标普全球高端消费品指数增强（QDII）美元现汇"
$q3Sheet.Range("D3").Value = "'1.85"
$q3Sheet.Range("E3").Value = "'93.04"
$q3Sheet.Range("F3").Value = "'9.69"
$q3Sheet.Range("G3").Value = "'0.1793"
$q3Sheet.Range("H3").Value = 1

$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("B4").Value = "'005676"
$q3Sheet.Range("C4").Value = "易方达标普全球高端消费品指数增强C（QDII）人民币"
$q3Sheet.Range("D4").Value = "'1.85"
$q3Sheet.Range("E4").Value = "'93.04"
$q3Sheet.Range("F4").Value = "'9.69"
$q3Sheet.Range("G4").Value = "'0.1793"
$q3Sheet.Range("H4").Value = 1

$q3Sheet.Range("A5").Value = 3
$q3Sheet.Range("B5").Value = "'513080"
$q3Sheet.Range("C5").Value = "华安法国CAC40ETF（QDII）"
$q3Sheet.Range("D5").Value = "'0.58"
$q3Sheet.Range("E5").Value = "'95.06"
$q3Sheet.Range("F5").Value = "'10.81"
$q3Sheet.Range("G5").Value = "'0.0627"
$q3Sheet.Range("H5").Value = 1

$q3Sheet.Range("A6").Value = 4
$q3Sheet.Range("B6").Value = "'006282"
$q3Sheet.Range("C6").Value = "上投摩根欧洲动力策略股票（QDII）"
$q3Sheet.Range("D6").Value = "'0.41"
$q3Sheet.Range("E6").Value = "'91.47"
$q3Sheet.Range("F6").Value = "'2.65"
$q3Sheet.Range("G6").Value = "'0.0109"
$q3Sheet.Range("H6").Value = 6

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a 2022-Q3 row right after the
#    header, pushing the existing rows down by one.
# ---------------------------------------------------------------------

# Copy the format of the last data row down onto the (currently empty) row 5
$totalSheet.Range("A4:D4").Copy()
$totalSheet.Range("A5:D5").PasteSpecial(-4122)

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 8
$totalSheet.Range("D5").Value = 2.5

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 8
$totalSheet.Range("D4").Value = 0.75

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 5
$totalSheet.Range("D3").Value = 0.58

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 0.61
